$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 411.61905
$ws.Range("I19").Value = 199.14285
$ws.Range("J19").Value = 517.8570999999999
$ws.Range("K19").Value = 199.14285
$ws.Range("L19").Value = 517.8570999999999
$ws.Range("M19").Value = -24.14285000000001
$ws.Range("N19").Value = -867.8570999999999
$ws.Range("H112").Value = 1303.8096
$ws.Range("I112").Value = 350
$ws.Range("J112").Value = 1780.7142
$ws.Range("K112").Value = 1050
$ws.Range("L112").Value = 5342.142599999999
$ws.Range("M112").Value = 58
$ws.Range("N112").Value = -7558.142599999999
$ws.Range("H135").Value = 295525.53
$ws.Range("I135").Value = 324080.56
$ws.Range("K135").Value = 2916725.04
$ws.Range("M135").Value = -2914190.04
$ws.Range("H141").Value = 2585.625
$ws.Range("I141").Value = 2859.1667
$ws.Range("J141").Value = 2494.4443
$ws.Range("K141").Value = 8577.500100000001
$ws.Range("L141").Value = 7483.3329
$ws.Range("M141").Value = -3397.500100000001
$ws.Range("N141").Value = -17843.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3905.81
$ws.Range("I32").Value = 3694.704
$ws.Range("J32").Value = 14250
$ws.Range("K32").Value = 3694.704
$ws.Range("L32").Value = 14250
$ws.Range("M32").Value = -3407.704
$ws.Range("N32").Value = -14824
$ws.Range("H138").Value = 54198.332
$ws.Range("J138").Value = 54198.332
$ws.Range("L138").Value = 54198.332
$ws.Range("N138").Value = -64478.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 23342.312
$ws.Range("I82").Value = 11038.8
$ws.Range("J82").Value = 28934.818
$ws.Range("K82").Value = 11038.8
$ws.Range("L82").Value = 28934.818
$ws.Range("M82").Value = -10655.8
$ws.Range("N82").Value = -29700.818
$ws.Range("H85").Value = 23342.312
$ws.Range("I85").Value = 11038.8
$ws.Range("J85").Value = 28934.818
$ws.Range("K85").Value = 11038.8
$ws.Range("L85").Value = 28934.818
$ws.Range("M85").Value = -9712.799999999999
$ws.Range("N85").Value = -31586.818
$ws.Range("H122").Value = 45331.25
$ws.Range("J122").Value = 45331.25
$ws.Range("L122").Value = 45331.25
$ws.Range("N122").Value = -55131.25
$ws.Range("H123").Value = 39269.332
$ws.Range("J123").Value = 39269.332
$ws.Range("L123").Value = 39269.332
$ws.Range("N123").Value = -49069.332
$ws.Range("H125").Value = 51451.668
$ws.Range("J125").Value = 51451.668
$ws.Range("L125").Value = 51451.668
$ws.Range("N125").Value = -61291.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 9100.5
$ws.Range("J50").Value = 9100.5
$ws.Range("L50").Value = 9100.5
$ws.Range("N50").Value = -10350.5
$ws.Range("H138").Value = 36640
$ws.Range("J138").Value = 36640
$ws.Range("L138").Value = 36640
$ws.Range("N138").Value = -46920
$ws.Range("H139").Value = 54998.5
$ws.Range("J139").Value = 54998.5
$ws.Range("L139").Value = 54998.5
$ws.Range("N139").Value = -65278.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 407.46667
$ws.Range("I12").Value = 67.09999999999999
$ws.Range("J12").Value = 577.65
$ws.Range("K12").Value = 201.3
$ws.Range("L12").Value = 1732.95
$ws.Range("M12").Value = -28.29999999999998
$ws.Range("N12").Value = -2078.95
$ws.Range("H14").Value = 491.7619
$ws.Range("I14").Value = 491.7619
$ws.Range("K14").Value = 1475.2857
$ws.Range("M14").Value = -1302.2857
$ws.Range("H114").Value = 609.2
$ws.Range("I114").Value = 817
$ws.Range("J114").Value = 297.5
$ws.Range("K114").Value = 2451
$ws.Range("L114").Value = 892.5
$ws.Range("M114").Value = 803
$ws.Range("N114").Value = -7400.5
$ws.Range("H117").Value = 3297.5
$ws.Range("I117").Value = 1500
$ws.Range("J117").Value = 3435.7693
$ws.Range("K117").Value = 4500
$ws.Range("L117").Value = 10307.3079
$ws.Range("N117").Value = -17191.3079
$ws.Range("M117").Value = -1058
$ws.Range("H121").Value = 580868
$ws.Range("I121").Value = 500
$ws.Range("J121").Value = 1027304.94
$ws.Range("K121").Value = 1500
$ws.Range("L121").Value = 3081914.82
$ws.Range("N121").Value = -3084534.82
$ws.Range("M121").Value = -190
$ws.Range("H129").Value = 64895.375
$ws.Range("I129").Value = 2115
$ws.Range("J129").Value = 73864
$ws.Range("K129").Value = 6345
$ws.Range("L129").Value = 221592
$ws.Range("M129").Value = -1345
$ws.Range("N129").Value = -231592
$ws.Range("H131").Value = 940.5
$ws.Range("I131").Value = 552.5
$ws.Range("J131").Value = 974.23914
$ws.Range("K131").Value = 1657.5
$ws.Range("L131").Value = 2922.71742
$ws.Range("M131").Value = 3382.5
$ws.Range("N131").Value = -13002.71742
$ws.Range("H132").Value = 442008.2
$ws.Range("I132").Value = 1013226.6
$ws.Range("J132").Value = 5194.1177
$ws.Range("K132").Value = 9119039.4
$ws.Range("L132").Value = 46747.0593
$ws.Range("M132").Value = -9116509.4
$ws.Range("N132").Value = -51807.0593

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 10000
$ws.Range("J52").Value = 10000
$ws.Range("L52").Value = 10000
$ws.Range("N52").Value = -10518
$ws.Range("H57").Value = 15798
$ws.Range("J57").Value = 17946.6
$ws.Range("L57").Value = 17946.6
$ws.Range("N57").Value = -19586.6
$ws.Range("H113").Value = 1807.0869
$ws.Range("I113").Value = 1657.0588
$ws.Range("J113").Value = 2232.1667
$ws.Range("K113").Value = 1657.0588
$ws.Range("L113").Value = 2232.1667
$ws.Range("M113").Value = 512.9412
$ws.Range("N113").Value = -6572.1667
$ws.Range("H136").Value = 18016.838
$ws.Range("J136").Value = 18016.838
$ws.Range("L136").Value = 54050.514
$ws.Range("N136").Value = -59150.514

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8401.6
$ws.Range("I81").Value = 21204.8
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 42409.6
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = -41348.6
$ws.Range("N81").Value = -6122
$ws.Range("H84").Value = 8401.6
$ws.Range("I84").Value = 21204.8
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 212048
$ws.Range("L84").Value = 20000
$ws.Range("M84").Value = -206744
$ws.Range("N84").Value = -30608
$ws.Range("H132").Value = 21846.1
$ws.Range("I132").Value = 2466.423
$ws.Range("J132").Value = 42840.75
$ws.Range("K132").Value = 7399.268999999999
$ws.Range("L132").Value = 128522.25
$ws.Range("M132").Value = -4869.268999999999
$ws.Range("N132").Value = -133582.25
